$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.034288629020576
$ws.Range("D2").Value = 1.036445064899648
$ws.Range("E2").Value = 1.04308523591621
$ws.Range("F2").Value = 1.053144062657496
$ws.Range("I2").Value = 1.036169870978543
$ws.Range("J2").Value = 1.039408692253994
$ws.Range("K2").Value = 1.039238786654378
$ws.Range("L2").Value = 1.045860096804635
$ws.Range("M2").Value = 1.055890843412746
$ws.Range("N2").Value = 1.040884771818838
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.035140055680802
$ws.Range("D3").Value = 1.037058555438295
$ws.Range("E3").Value = 1.043915656225188
$ws.Range("F3").Value = 1.054223695634863
$ws.Range("I3").Value = 1.036342862173903
$ws.Range("J3").Value = 1.03990368429834
$ws.Range("K3").Value = 1.039662400553162
$ws.Range("L3").Value = 1.046501429441061
$ws.Range("M3").Value = 1.056782779266757
$ws.Range("N3").Value = 1.041380466808664
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.03569147886536
$ws.Range("D4").Value = 1.03745589218931
$ws.Range("E4").Value = 1.044453875785342
$ws.Range("F4").Value = 1.054923705310342
$ws.Range("I4").Value = 1.036453829170416
$ws.Range("J4").Value = 1.040223814536006
$ws.Range("K4").Value = 1.039936180652066
$ws.Range("L4").Value = 1.046916647824395
$ws.Range("M4").Value = 1.057360737396848
$ws.Range("N4").Value = 1.041701051667987
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.035923413868025
$ws.Range("D5").Value = 1.03762301893842
$ws.Range("E5").Value = 1.04468035298415
$ws.Range("F5").Value = 1.055218326777118
$ws.Range("I5").Value = 1.036500246909858
$ws.Range("J5").Value = 1.040358357255754
$ws.Range("K5").Value = 1.040051198590457
$ws.Range("L5").Value = 1.047091260304403
$ws.Range("M5").Value = 1.057603905077824
$ws.Range("N5").Value = 1.041835785453829
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.035962363578666
$ws.Range("D6").Value = 1.037651085274752
$ws.Range("E6").Value = 1.044718391756086
$ws.Range("F6").Value = 1.05526781476416
$ws.Range("I6").Value = 1.036508026984593
$ws.Range("J6").Value = 1.040380945191595
$ws.Range("K6").Value = 1.040070505932114
$ws.Range("L6").Value = 1.047120581679938
$ws.Range("M6").Value = 1.057644745345851
$ws.Range("N6").Value = 1.041858405467129
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.035694577535401
$ws.Range("D7").Value = 1.037458125007669
$ws.Range("E7").Value = 1.044456901162648
$ws.Range("F7").Value = 1.054927640732724
$ws.Range("I7").Value = 1.036454450321952
$ws.Range("J7").Value = 1.040225612461061
$ws.Range("K7").Value = 1.03993771783999
$ws.Range("L7").Value = 1.046918980791549
$ws.Range("M7").Value = 1.057363985854269
$ws.Range("N7").Value = 1.041702852146301
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.034576269696624
$ws.Range("D8").Value = 1.036652319896812
$ws.Range("E8").Value = 1.043365696815078
$ws.Range("F8").Value = 1.053508635894619
$ws.Range("I8").Value = 1.036228534623031
$ws.Range("J8").Value = 1.039576010175473
$ws.Range("K8").Value = 1.039382015899002
$ws.Range("L8").Value = 1.046076789160914
$ws.Range("M8").Value = 1.05619210759146
$ws.Range("N8").Value = 1.041052327350957
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.032609506711689
$ws.Range("D9").Value = 1.035235270933924
$ws.Range("E9").Value = 1.041449676515718
$ws.Range("F9").Value = 1.05101906430851
$ws.Range("I9").Value = 1.035823040615376
$ws.Range("J9").Value = 1.038430129425247
$ws.Range("K9").Value = 1.038400345375856
$ws.Range("L9").Value = 1.04459457863442
$ws.Range("M9").Value = 1.054133413975326
$ws.Range("N9").Value = 1.039904819318632
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.031300993699358
$ws.Range("D10").Value = 1.034292602869719
$ws.Range("E10").Value = 1.040177003667018
$ws.Range("F10").Value = 1.049366755318007
$ws.Range("I10").Value = 1.035547770299876
$ws.Range("J10").Value = 1.037665467121572
$ws.Range("K10").Value = 1.037744312885033
$ws.Range("L10").Value = 1.043607744383338
$ws.Range("M10").Value = 1.052765266132811
$ws.Range("N10").Value = 1.039139071106777
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.030735041617805
$ws.Range("D11").Value = 1.033884918604616
$ws.Range("E11").Value = 1.039627049302413
$ws.Range("F11").Value = 1.048653060239072
$ws.Range("I11").Value = 1.035427410718434
$ws.Range("J11").Value = 1.037334197368416
$ws.Range("K11").Value = 1.037459880869673
$ws.Range("L11").Value = 1.043180758657839
$ws.Range("M11").Value = 1.052173881994815
$ws.Range("N11").Value = 1.038807330912575
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.030524919720326
$ws.Range("D12").Value = 1.033733562801622
$ws.Range("E12").Value = 1.039422941713218
$ws.Range("F12").Value = 1.04838822867758
$ws.Range("I12").Value = 1.035382529254765
$ws.Range("J12").Value = 1.037211125167753
$ws.Range("K12").Value = 1.037354176286983
$ws.Range("L12").Value = 1.043022206386995
$ws.Range("M12").Value = 1.0519543720122
$ws.Range("N12").Value = 1.03868408393527
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.030569987128535
$ws.Range("D13").Value = 1.033766025657703
$ws.Range("E13").Value = 1.03946671576974
$ws.Range("F13").Value = 1.048445023862743
$ws.Range("I13").Value = 1.035392164370777
$ws.Range("J13").Value = 1.037237525635391
$ws.Range("K13").Value = 1.037376852699369
$ws.Range("L13").Value = 1.043056214140482
$ws.Range("M13").Value = 1.052001450550506
$ws.Range("N13").Value = 1.0387105218946
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.030717670871174
$ws.Range("D14").Value = 1.033872405921386
$ws.Range("E14").Value = 1.039610174229986
$ws.Range("F14").Value = 1.048631163734111
$ws.Range("I14").Value = 1.035423704360177
$ws.Range("J14").Value = 1.037324024662585
$ws.Range("K14").Value = 1.037451144388637
$ws.Range("L14").Value = 1.043167651666487
$ws.Range("M14").Value = 1.05215573403314
$ws.Range("N14").Value = 1.038797143760334
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.030808676712904
$ws.Range("D15").Value = 1.033937960468658
$ws.Range("E15").Value = 1.039698586286034
$ws.Range("F15").Value = 1.04874588601423
$ws.Range("I15").Value = 1.035443114075601
$ws.Range("J15").Value = 1.037377316437647
$ws.Range("K15").Value = 1.037496910844095
$ws.Range("L15").Value = 1.043236318563573
$ws.Range("M15").Value = 1.052250813942678
$ws.Range("N15").Value = 1.038850511215831
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.031338567730758
$ws.Range("D16").Value = 1.03431967014285
$ws.Range("E16").Value = 1.040213526078895
$ws.Range("F16").Value = 1.049414158220597
$ws.Range("I16").Value = 1.035555733658008
$ws.Range("J16").Value = 1.037687448998964
$ws.Range("K16").Value = 1.03776318208864
$ws.Range("L16").Value = 1.043636088861067
$ws.Range("M16").Value = 1.052804536217834
$ws.Range("N16").Value = 1.039161084200957
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.031671127382851
$ws.Range("D17").Value = 1.034559240776477
$ws.Range("E17").Value = 1.040536835471134
$ws.Range("F17").Value = 1.049833821344713
$ws.Range("I17").Value = 1.035626065259508
$ws.Range("J17").Value = 1.037881943288068
$ws.Range("K17").Value = 1.037930109868924
$ws.Range("L17").Value = 1.043886940702875
$ws.Range("M17").Value = 1.05315214898059
$ws.Range("N17").Value = 1.03935585469426
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.031865165720871
$ws.Range("D18").Value = 1.034699026143082
$ws.Range("E18").Value = 1.040725524467335
$ws.Range("F18").Value = 1.050078774013565
$ws.Range("I18").Value = 1.035666975978831
$ws.Range("J18").Value = 1.037995372435052
$ws.Range("K18").Value = 1.038027440670912
$ws.Range("L18").Value = 1.044033289137217
$ws.Range("M18").Value = 1.053355005201323
$ws.Range("N18").Value = 1.039469444923643
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.03193133824139
$ws.Range("D19").Value = 1.034746697411995
$ws.Range("E19").Value = 1.040789880829297
$ws.Range("F19").Value = 1.050162325384446
$ws.Range("I19").Value = 1.035680906374969
$ws.Range("J19").Value = 1.038034046082874
$ws.Range("K19").Value = 1.038060621931189
$ws.Range("L19").Value = 1.044083195360147
$ws.Range("M19").Value = 1.053424190782867
$ws.Range("N19").Value = 1.039508173492481
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.031635440465136
$ws.Range("D20").Value = 1.034533532151998
$ws.Range("E20").Value = 1.040502136241417
$ws.Range("F20").Value = 1.049788777842025
$ws.Range("I20").Value = 1.035618530977936
$ws.Range("J20").Value = 1.037861077547193
$ws.Range("K20").Value = 1.037912203746434
$ws.Range("L20").Value = 1.043860023480209
$ws.Range("M20").Value = 1.053114843112236
$ws.Range("N20").Value = 1.03933495932164
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.030674178976394
$ws.Range("D21").Value = 1.033841077466647
$ws.Range("E21").Value = 1.039567924585792
$ws.Range("F21").Value = 1.04857634281232
$ws.Range("I21").Value = 1.035414421433348
$ws.Range("J21").Value = 1.037298553505231
$ws.Range("K21").Value = 1.037429268816304
$ws.Range("L21").Value = 1.043134834724016
$ws.Range("M21").Value = 1.052110297065084
$ws.Range("N21").Value = 1.038771636431016
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.030070362842398
$ws.Range("D22").Value = 1.033406145866963
$ws.Range("E22").Value = 1.038981532640165
$ws.Range("F22").Value = 1.047815579718679
$ws.Range("I22").Value = 1.035285080022304
$ws.Range("J22").Value = 1.036944734768513
$ws.Range("K22").Value = 1.037125317321443
$ws.Range("L22").Value = 1.042679165228955
$ws.Range("M22").Value = 1.051479604694148
$ws.Range("N22").Value = 1.038417315231108
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.030390402992304
$ws.Range("D23").Value = 1.03363666889912
$ws.Range("E23").Value = 1.039292296340733
$ws.Range("F23").Value = 1.048218727939622
$ws.Range("I23").Value = 1.035353741859484
$ws.Range("J23").Value = 1.037132313457194
$ws.Range("K23").Value = 1.03728647689721
$ws.Range("L23").Value = 1.042920696803623
$ws.Range("M23").Value = 1.051813860398973
$ws.Range("N23").Value = 1.038605160303042
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.031651565655946
$ws.Range("D24").Value = 1.034545148625073
$ws.Range("E24").Value = 1.040517814995877
$ws.Range("F24").Value = 1.049809130536033
$ws.Range("I24").Value = 1.035621935743182
$ws.Range("J24").Value = 1.037870505927145
$ws.Range("K24").Value = 1.037920294861769
$ws.Range("L24").Value = 1.043872186119452
$ws.Range("M24").Value = 1.053131699720773
$ws.Range("N24").Value = 1.039344401090973
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.03311749882809
$ws.Range("D25").Value = 1.035601260849427
$ws.Range("E25").Value = 1.041944196300682
$ws.Range("F25").Value = 1.051661378820482
$ws.Range("I25").Value = 1.035928744023307
$ws.Range("J25").Value = 1.038726502505563
$ws.Range("K25").Value = 1.038654414844707
$ws.Range("L25").Value = 1.044977540497086
$ws.Range("M25").Value = 1.054664880661165
$ws.Range("N25").Value = 1.040201613282716
